$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear O4 (header "Internal Assignment")
$ws.Range("O4").ClearContents()

# Clear O5:O7 (values under "Internal Assignment")
$ws.Range("O5:O7").ClearContents()

# Change A5 value from $NAME to $$NAME
$ws.Range("A5").Value = "`$`$NAME"

# Update selection to O4:O7 with active cell O4
$ws.Range("O4:O7").Select()
